$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 58201.37
$ws.Range("J17").Value = 58201.37
$ws.Range("L17").Value = 174604.11
$ws.Range("N17").Value = -174940.11
$ws.Range("H19").Value = 2375.923
$ws.Range("I19").Value = 2224.5
$ws.Range("K19").Value = 2224.5
$ws.Range("M19").Value = -2049.5
$ws.Range("H28").Value = 33369.805
$ws.Range("I28").Value = 50876.35
$ws.Range("J28").Value = 1539.7273
$ws.Range("K28").Value = 50876.35
$ws.Range("L28").Value = 1539.7273
$ws.Range("M28").Value = -50391.35
$ws.Range("N28").Value = -2509.7273
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H41").Value = 20001760
$ws.Range("I41").Value = 596.13336
$ws.Range("K41").Value = 596.13336
$ws.Range("M41").Value = -156.13336
$ws.Range("H43").Value = 1950.5
$ws.Range("I43").Value = 1950.5
$ws.Range("K43").Value = 1950.5
$ws.Range("M43").Value = -1881.5
$ws.Range("H53").Value = 13334096
$ws.Range("I53").Value = 55556476
$ws.Range("J53").Value = 712.3684
$ws.Range("K53").Value = 55556476
$ws.Range("L53").Value = 712.3684
$ws.Range("M53").Value = -55555839
$ws.Range("N53").Value = -1986.3684
$ws.Range("H80").Value = 5114.4614
$ws.Range("I80").Value = 7499.1333
$ws.Range("J80").Value = 1862.6364
$ws.Range("K80").Value = 22497.3999
$ws.Range("L80").Value = 5587.9092
$ws.Range("M80").Value = -21499.3999
$ws.Range("N80").Value = -7583.9092
$ws.Range("H83").Value = 5114.4614
$ws.Range("I83").Value = 7499.1333
$ws.Range("J83").Value = 1862.6364
$ws.Range("K83").Value = 67492.1997
$ws.Range("L83").Value = 16763.7276
$ws.Range("M83").Value = -62500.1997
$ws.Range("N83").Value = -26747.7276
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H98").Value = 2473.1738
$ws.Range("I98").Value = 1125.8125
$ws.Range("J98").Value = 5552.857
$ws.Range("K98").Value = 1125.8125
$ws.Range("L98").Value = 5552.857
$ws.Range("M98").Value = 372.1875
$ws.Range("N98").Value = -8548.857
$ws.Range("H122").Value = 2473.1738
$ws.Range("I122").Value = 1125.8125
$ws.Range("J122").Value = 5552.857
$ws.Range("K122").Value = 3377.4375
$ws.Range("L122").Value = 16658.571
$ws.Range("M122").Value = -927.4375
$ws.Range("N122").Value = -21558.571
$ws.Range("H125").Value = 7815.643
$ws.Range("I125").Value = 7449.1
$ws.Range("K125").Value = 67041.90000000001
$ws.Range("M125").Value = -64581.90000000001
$ws.Range("H132").Value = 7062.44
$ws.Range("I132").Value = 8898.429
$ws.Range("K132").Value = 26695.287
$ws.Range("M132").Value = -24165.287
$ws.Range("H137").Value = 2706.975
$ws.Range("I137").Value = 1779.32
$ws.Range("K137").Value = 5337.96
$ws.Range("M137").Value = -2787.96
$ws.Range("H138").Value = 4265.4575
$ws.Range("I138").Value = 940.5789
$ws.Range("J138").Value = 5844.775
$ws.Range("K138").Value = 2821.7367
$ws.Range("L138").Value = 17534.325
$ws.Range("M138").Value = 2318.2633
$ws.Range("N138").Value = -27814.325

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10328.827
$ws.Range("I74").Value = 11585.917
$ws.Range("J74").Value = 4294.8
$ws.Range("K74").Value = 11585.917
$ws.Range("L74").Value = 4294.8
$ws.Range("M74").Value = -10711.917
$ws.Range("N74").Value = -6042.8
$ws.Range("H77").Value = 10328.827
$ws.Range("I77").Value = 11585.917
$ws.Range("J77").Value = 4294.8
$ws.Range("K77").Value = 57929.585
$ws.Range("L77").Value = 21474
$ws.Range("M77").Value = -53561.585
$ws.Range("N77").Value = -30210
$ws.Range("H126").Value = 9999.5
$ws.Range("I126").Value = 9999.5
$ws.Range("K126").Value = 29998.5
$ws.Range("M126").Value = -27528.5
$ws.Range("H132").Value = 7015.6313
$ws.Range("I132").Value = 5142.44
$ws.Range("K132").Value = 15427.32
$ws.Range("M132").Value = -12897.32
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360
$ws.Range("H141").Value = 65000
$ws.Range("J141").Value = 65000
$ws.Range("L141").Value = 65000
$ws.Range("N141").Value = -75360

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 215
$ws.Range("I22").Value = 240
$ws.Range("J22").Value = 140
$ws.Range("K22").Value = 240
$ws.Range("L22").Value = 140
$ws.Range("M22").Value = -67
$ws.Range("N22").Value = -486
$ws.Range("H94").Value = 956.4
$ws.Range("I94").Value = 754.4286
$ws.Range("J94").Value = 1427.6666
$ws.Range("K94").Value = 754.4286
$ws.Range("L94").Value = 1427.6666
$ws.Range("M94").Value = -303.4286
$ws.Range("N94").Value = -2329.6666
$ws.Range("H134").Value = 4610.2354
$ws.Range("I134").Value = 3197.8333
$ws.Range("K134").Value = 9593.499899999999
$ws.Range("M134").Value = -7058.499899999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4803.0586
$ws.Range("I99").Value = 2731.625
$ws.Range("J99").Value = 6644.3335
$ws.Range("K99").Value = 2731.625
$ws.Range("L99").Value = 6644.3335
$ws.Range("M99").Value = -1233.625
$ws.Range("N99").Value = -9640.333500000001
$ws.Range("H122").Value = 1931.1212
$ws.Range("J122").Value = 3851.4443
$ws.Range("L122").Value = 11554.3329
$ws.Range("N122").Value = -16454.3329
$ws.Range("H126").Value = 4803.0586
$ws.Range("I126").Value = 2731.625
$ws.Range("J126").Value = 6644.3335
$ws.Range("K126").Value = 8194.875
$ws.Range("L126").Value = 19933.0005
$ws.Range("M126").Value = -5724.875
$ws.Range("N126").Value = -24873.0005
$ws.Range("H132").Value = 3153.6667
$ws.Range("I132").Value = 2459.7932
$ws.Range("K132").Value = 7379.3796
$ws.Range("M132").Value = -4849.3796

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2200
$ws.Range("I80").Value = 2200
$ws.Range("K80").Value = 6600
$ws.Range("M80").Value = -5664
$ws.Range("H83").Value = 2200
$ws.Range("I83").Value = 2200
$ws.Range("K83").Value = 19800
$ws.Range("M83").Value = -15120
$ws.Range("H92").Value = 598
$ws.Range("I92").Value = 166.33333
$ws.Range("J92").Value = 759.875
$ws.Range("K92").Value = 498.99999
$ws.Range("L92").Value = 2279.625
$ws.Range("M92").Value = 749.00001
$ws.Range("N92").Value = -4775.625
$ws.Range("H132").Value = 2859.88
$ws.Range("I132").Value = 1188.1765
$ws.Range("J132").Value = 6412.25
$ws.Range("K132").Value = 10693.5885
$ws.Range("L132").Value = 57710.25
$ws.Range("M132").Value = -8163.5885
$ws.Range("N132").Value = -62770.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 134.7
$ws.Range("I2").Value = 134.7
$ws.Range("K2").Value = 134.7
$ws.Range("M2").Value = -21.69999999999999
$ws.Range("H20").Value = 32500
$ws.Range("J20").Value = 32500
$ws.Range("L20").Value = 32500
$ws.Range("N20").Value = -32990
$ws.Range("H43").Value = 26600
$ws.Range("J43").Value = 34900
$ws.Range("L43").Value = 34900
$ws.Range("N43").Value = -35202
$ws.Range("H80").Value = 6667837
$ws.Range("I80").Value = 5001752.5
$ws.Range("J80").Value = 10000006
$ws.Range("K80").Value = 5001752.5
$ws.Range("L80").Value = 10000006
$ws.Range("M80").Value = -5000754.5
$ws.Range("N80").Value = -10002002
$ws.Range("H83").Value = 6667837
$ws.Range("I83").Value = 5001752.5
$ws.Range("J83").Value = 10000006
$ws.Range("K83").Value = 25008762.5
$ws.Range("L83").Value = 50000030
$ws.Range("M83").Value = -25003770.5
$ws.Range("N83").Value = -50010014
$ws.Range("H122").Value = 6009.1724
$ws.Range("I122").Value = 5078.778
$ws.Range("K122").Value = 15236.334
$ws.Range("M122").Value = -12786.334
$ws.Range("H131").Value = 35122.125
$ws.Range("J131").Value = 34425.285
$ws.Range("L131").Value = 34425.285
$ws.Range("N131").Value = -44505.285

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 460646.7
$ws.Range("I7").Value = 671482.6
$ws.Range("K7").Value = 671482.6
$ws.Range("M7").Value = -671370.6
$ws.Range("H82").Value = 2450.7856
$ws.Range("I82").Value = 2554.1428
$ws.Range("J82").Value = 2347.4285
$ws.Range("K82").Value = 2554.1428
$ws.Range("L82").Value = 2347.4285
$ws.Range("M82").Value = -2193.1428
$ws.Range("N82").Value = -3069.4285
$ws.Range("H85").Value = 2450.7856
$ws.Range("I85").Value = 2554.1428
$ws.Range("J85").Value = 2347.4285
$ws.Range("K85").Value = 2554.1428
$ws.Range("L85").Value = 2347.4285
$ws.Range("M85").Value = -1306.1428
$ws.Range("N85").Value = -4843.4285
$ws.Range("H122").Value = 681550.75
$ws.Range("I122").Value = 558317.9399999999
$ws.Range("K122").Value = 1674953.82
$ws.Range("M122").Value = -1672503.82
$ws.Range("H126").Value = 460646.7
$ws.Range("I126").Value = 671482.6
$ws.Range("K126").Value = 2014447.8
$ws.Range("M126").Value = -2011977.8
$ws.Range("H132").Value = 5899.8887
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8790.1
$ws.Range("I41").Value = 8788.5
$ws.Range("J41").Value = 8790.5
$ws.Range("K41").Value = 8788.5
$ws.Range("L41").Value = 8790.5
$ws.Range("M41").Value = -8398.5
$ws.Range("N41").Value = -9570.5
$ws.Range("H126").Value = 2801.9565
$ws.Range("I126").Value = 2365.7896
$ws.Range("K126").Value = 7097.3688
$ws.Range("M126").Value = -4627.3688
